$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply only the cells whose value actually changes, cell-by-cell,
# matching the target diff exactly (rows 8-21 reshuffled, row 22 added).

# Row 2
$ws.Range("D2").Value = "5178, 5582, 6488"
# Row 8
$ws.Range("A8").Value = "98, 130, 1073, SF"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "'130"
$ws.Range("D8").Value = "'4415"
# Row 9
$ws.Range("A9").Value = "423, 1073, 1105, SF, SF"
$ws.Range("C9").Value = "'1105"
$ws.Range("D9").Value = "'4994"
# Row 11
$ws.Range("A11").Value = "130, 748, 1073, 1105"
$ws.Range("C11").Value = "'1105"
$ws.Range("D11").Value = "'5399"
# Row 12
$ws.Range("A12").Value = "423, 748, 780, 1073"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "'780"
$ws.Range("D12").Value = "'5677"
# Row 13
$ws.Range("A13").Value = "423, 748, 1073, 1105"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = "1105, 1105"
$ws.Range("D13").Value = "5331, 5433"
# Row 14
$ws.Range("A14").Value = "130, 423, 748, 1073"
$ws.Range("D14").Value = "'5886"
# Row 15
$ws.Range("A15").Value = "98, 130, 455, 748, 1073"
$ws.Range("D15").Value = "'6016"
# Row 16
$ws.Range("A16").Value = "130, 423, 748, 1073, SF"
$ws.Range("D16").Value = "'6202"
# Row 17
$ws.Range("A17").Value = "130, 130, 423"
$ws.Range("D17").Value = "'6426"
# Row 18
$ws.Range("A18").Value = "130, 455, 780, 1105"
$ws.Range("C18").Value = "'130"
$ws.Range("D18").Value = "'6561"
# Row 19
$ws.Range("A19").Value = "98, 780, 780, 1105"
$ws.Range("C19").Value = "'780"
$ws.Range("D19").Value = "'6359"
# Row 20
$ws.Range("A20").Value = "98, 130, 130, 455, 780"
$ws.Range("D20").Value = "'6727"
# Row 21
$ws.Range("A21").Value = "130, 130, 423, 1073"
$ws.Range("C21").Value = "'130"
$ws.Range("D21").Value = "'6670"
# Row 22
$ws.Range("A22").Value = "98, 98, 455, 455"
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = "'455"
$ws.Range("D22").Value = "'6748"

# These cells hold numeric-looking text (e.g. "130", "4415"); a leading
# apostrophe was used above to force Excel to store them as text rather
# than coercing them into numbers. That leaves a quote-prefix number
# format on the cell, so clear the formatting back off (per cell, since
# ClearFormats on a multi-area Union range only affects the first area)
# while the stored value stays text.
$ws.Range("C8").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("C9").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("C11").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("C12").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("C18").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("C19").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("C21").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("C22").ClearFormats()
$ws.Range("D22").ClearFormats()

